$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New summary row (row 54): label + MAX/MIN formulas over the data rows (2-51)
$ws.Range("A54").Value = "min/max"

$ws.Range("B54").Formula = "=MAX(B1:B51)"
$ws.Range("C54").Formula = "=MIN(C1:C51)"
$ws.Range("F54").Formula = "=MAX(F1:F51)"
$ws.Range("H54").Formula = "=MIN(H1:H51)"

# Column widths tweaked (auto-fit like values) for columns B-F and H
$ws.Range("B1").ColumnWidth = 14.85546875
$ws.Range("C1").ColumnWidth = 17.7109375
$ws.Range("D1").ColumnWidth = 19
$ws.Range("E1").ColumnWidth = 18.28515625
$ws.Range("F1").ColumnWidth = 18.5703125
$ws.Range("H1").ColumnWidth = 28

# Scroll/selection state matches the saved view in the workbook
$ws.Application.ActiveWindow.ScrollRow = 37
$ws.Range("F54").Select()
